# AdminSection Module datepickerUtil: refresh the ECM test user fixture
# (Ecm04 -> Ecm06 / testuser001 -> testuser006) on the TestData sheet,
# and move the saved selection to E32.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestData")

# ECM_FirstName (row 32)
$ws.Range("B32").Value = "Ecm06"

# ECM_UserName (row 34)
$ws.Range("B34").Value = "Ecm06_testuser"

# ECM_EmailId (row 35)
$ws.Range("B35").Value = "testuser006@test.com"

# EditUserName (row 41) - mirrors ECM_UserName
$ws.Range("B41").Value = "Ecm06_testuser"

# Update the saved selection/active cell for the sheet view
$ws.Range("E32").Select()
